$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 830.04
$ws.Range("C3").Value = 841.34
$ws.Range("C4").Value = 807.87
$ws.Range("C5").Value = 813.52
$ws.Range("C6").Value = 813.52
